$wb = $excel.ActiveWorkbook

# Hunk 0: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 975.1081
$ws.Range("J17").Value = 975.1081
$ws.Range("L17").Value = 2925.3243
$ws.Range("N17").Value = -3261.3243

# Hunk 1: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2572.0908
$ws.Range("J62").Value = 3350.5
$ws.Range("L62").Value = 3350.5
$ws.Range("N62").Value = -4598.5

# Hunk 2: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 2572.0908
$ws.Range("J65").Value = 3350.5
$ws.Range("L65").Value = 16752.5
$ws.Range("N65").Value = -22992.5

# Hunk 3: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4912.1333
$ws.Range("I76").Value = 4355.5713
$ws.Range("J76").Value = 5399.125
$ws.Range("K76").Value = 4355.5713
$ws.Range("L76").Value = 5399.125
$ws.Range("M76").Value = -4040.5713
$ws.Range("N76").Value = -6029.125

# Hunk 4: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4912.1333
$ws.Range("I79").Value = 4355.5713
$ws.Range("J79").Value = 5399.125
$ws.Range("K79").Value = 4355.5713
$ws.Range("L79").Value = 5399.125
$ws.Range("M79").Value = -3263.5713
$ws.Range("N79").Value = -7583.125

# Hunk 5: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4330.6665
$ws.Range("I86").Value = 4330.6665
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4330.6665
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3207.6665
$ws.Range("N86").Value = $null

# Hunk 6: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 4330.6665
$ws.Range("I89").Value = 4330.6665
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 21653.3325
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -16037.3325
$ws.Range("N89").Value = $null

# Hunk 7: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1450.9584
$ws.Range("I112").Value = 1213.3334
$ws.Range("J112").Value = 1847
$ws.Range("K112").Value = 3640.0002
$ws.Range("L112").Value = 5541
$ws.Range("M112").Value = -2532.0002
$ws.Range("N112").Value = -7757

# Hunk 8: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1002177.4
$ws.Range("I132").Value = 2254.1428
$ws.Range("J132").Value = 3335331.8
$ws.Range("K132").Value = 6762.428400000001
$ws.Range("L132").Value = 10005995.4
$ws.Range("M132").Value = -4232.428400000001
$ws.Range("N132").Value = -10011055.4

# Hunk 9: sheet ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2624.0396
$ws.Range("J138").Value = 2155.4915
$ws.Range("L138").Value = 6466.4745
$ws.Range("N138").Value = -16746.4745

# Hunk 10: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2521.875
$ws.Range("I45").Value = 955.3333
$ws.Range("K45").Value = 955.3333
$ws.Range("M45").Value = -578.3333

# Hunk 11: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1664.909
$ws.Range("I88").Value = 2289.3333
$ws.Range("K88").Value = 2289.3333
$ws.Range("M88").Value = -1883.3333

# Hunk 12: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1664.909
$ws.Range("I91").Value = 2289.3333
$ws.Range("K91").Value = 2289.3333
$ws.Range("M91").Value = -885.3332999999998

# Hunk 13: sheet ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1987.6154
$ws.Range("I132").Value = 1987.6154
$ws.Range("K132").Value = 5962.8462
$ws.Range("M132").Value = -3432.8462

# Hunk 14: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2738.652
$ws.Range("I20").Value = 1494.55
$ws.Range("K20").Value = 1494.55
$ws.Range("M20").Value = -1247.55

# Hunk 15: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null

# Hunk 16: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null

# Hunk 17: sheet BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3369.182
$ws.Range("I105").Value = 3152.2856
$ws.Range("J105").Value = 3748.75
$ws.Range("K105").Value = 3152.2856
$ws.Range("L105").Value = 3748.75
$ws.Range("M105").Value = -1405.2856
$ws.Range("N105").Value = -7242.75

# Hunk 18: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null

# Hunk 19: sheet CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null

# Hunk 20: sheet CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 267.76923
$ws.Range("I12").Value = 311.30768
$ws.Range("J12").Value = 224.23077
$ws.Range("K12").Value = 933.92304
$ws.Range("L12").Value = 672.69231
$ws.Range("M12").Value = -760.92304
$ws.Range("N12").Value = -1018.69231

# Hunk 21: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I70").Value = 9007
$ws.Range("J70").Value = 8846.25
$ws.Range("K70").Value = 9007
$ws.Range("L70").Value = 8846.25
$ws.Range("M70").Value = -8737
$ws.Range("N70").Value = -9386.25

# Hunk 22: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I73").Value = 9007
$ws.Range("J73").Value = 8846.25
$ws.Range("K73").Value = 9007
$ws.Range("L73").Value = 8846.25
$ws.Range("M73").Value = -8071
$ws.Range("N73").Value = -10718.25

# Hunk 23: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 12704.765
$ws.Range("I80").Value = 4634.625
$ws.Range("J80").Value = 19878.223
$ws.Range("K80").Value = 4634.625
$ws.Range("L80").Value = 19878.223
$ws.Range("M80").Value = -3636.625
$ws.Range("N80").Value = -21874.223

# Hunk 24: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 12704.765
$ws.Range("I83").Value = 4634.625
$ws.Range("J83").Value = 19878.223
$ws.Range("K83").Value = 23173.125
$ws.Range("L83").Value = 99391.11500000001
$ws.Range("M83").Value = -18181.125
$ws.Range("N83").Value = -109375.115

# Hunk 25: sheet GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2886.818
$ws.Range("I126").Value = 2500.4285
$ws.Range("J126").Value = 3563
$ws.Range("K126").Value = 7501.2855
$ws.Range("L126").Value = 10689
$ws.Range("M126").Value = -5031.2855
$ws.Range("N126").Value = -15629

# Hunk 26: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2511.8096
$ws.Range("I7").Value = 1431
$ws.Range("J7").Value = 3052.2144
$ws.Range("K7").Value = 1431
$ws.Range("L7").Value = 3052.2144
$ws.Range("M7").Value = -1319
$ws.Range("N7").Value = -3276.2144

# Hunk 27: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3135.75
$ws.Range("I40").Value = 2993.9167
$ws.Range("J40").Value = 3348.5
$ws.Range("K40").Value = 2993.9167
$ws.Range("L40").Value = 3348.5
$ws.Range("M40").Value = -2857.9167
$ws.Range("N40").Value = -3620.5

# Hunk 28: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3530.4614
$ws.Range("J46").Value = 3799.9
$ws.Range("L46").Value = 3799.9
$ws.Range("N46").Value = -4175.9

# Hunk 29: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3761.7812
$ws.Range("I122").Value = 4042.9
$ws.Range("K122").Value = 12128.7
$ws.Range("M122").Value = -9678.700000000001

# Hunk 30: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2511.8096
$ws.Range("I126").Value = 1431
$ws.Range("J126").Value = 3052.2144
$ws.Range("K126").Value = 4293
$ws.Range("L126").Value = 9156.643199999999
$ws.Range("M126").Value = -1823
$ws.Range("N126").Value = -14096.6432

# Hunk 31: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4440.7144
$ws.Range("J132").Value = 5700
$ws.Range("L132").Value = 17100
$ws.Range("N132").Value = -22160

# Hunk 32: sheet LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1962.3334
$ws.Range("I136").Value = 1962.3334
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5887.0002
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -3337.0002
$ws.Range("N136").Value = $null

# Hunk 33: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1880.125
$ws.Range("J81").Value = 850
$ws.Range("L81").Value = 1700
$ws.Range("N81").Value = -3822

# Hunk 34: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1880.125
$ws.Range("J84").Value = 850
$ws.Range("L84").Value = 8500
$ws.Range("N84").Value = -19108

# Hunk 35: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7643.1
$ws.Range("J122").Value = 4596.8335
$ws.Range("L122").Value = 13790.5005
$ws.Range("N122").Value = -18690.5005

# Hunk 36: sheet WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1545.4375
$ws.Range("I132").Value = 1545.4375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4636.3125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2106.3125
$ws.Range("N132").Value = $null
